$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B78 was stored as text "3" - correct it to a numeric 3 (matches the rest of column B)
$ws.Range("B78").Value = 3

# Append a new annotation row (row 79) for Ruilin
$ws.Range("A79").Value = "Ruilin"

# B79 keeps its politeness score as text "4" (like the header row), not a number
$ws.Range("B79").NumberFormat = "@"
$ws.Range("B79").Value = "4"
$ws.Range("B79").Style = "Normal"

$ws.Range("C79").Value = "The paper's contributions are significant."
$ws.Range("D79").Value = "APC"
$ws.Range("E79").Value = "OTH"
$ws.Range("F79").Value = "b9d28a3e-28bc-41b5-b6f1-68624390902f"
$ws.Range("G79").Value = "BJlrSmbAZ_annotated.xlsx"
$ws.Range("H79").Value = "The paper's contributions are significant."
